# Update model parameter estimates and covariance matrices
# with refreshed example output for HR and survival plots.

$wb = $excel.ActiveWorkbook

# --- sheet: weibull ---
$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -2.14218267845162
$ws.Range("C2").Value = 0.273912315481934
$ws.Range("B3").Value = 0.1969390943671
$ws.Range("C3").Value = 0.161633137077179

# --- sheet: lognormal ---
$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 1.49818030642701
$ws.Range("C2").Value = 0.331222567339514
$ws.Range("B3").Value = -1.01418267661033
$ws.Range("C3").Value = 0.147453826849644

# --- sheet: llogis ---
$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -1.51780137431937
$ws.Range("C2").Value = 0.174083053785454
$ws.Range("B3").Value = 0.580710717024211
$ws.Range("C3").Value = 0.149819611981135

# --- sheet: gompertz ---
$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -2.02381715538323
$ws.Range("C2").Value = 0.255778640267797
$ws.Range("B3").Value = 0.0255597996111659
$ws.Range("C3").Value = 0.0404981984260776

# --- sheet: weibull cov ---
$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.0750279565726748
$ws.Range("B2").Value = -0.0311812098035073
$ws.Range("A3").Value = -0.0311812098035073
$ws.Range("B3").Value = 0.0261252710014103

# --- sheet: lognormal cov ---
$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.109708389114979
$ws.Range("B2").Value = -0.0418133861145634
$ws.Range("A3").Value = -0.0418133861145634
$ws.Range("B3").Value = 0.0217426310526049

# --- sheet: llogis cov ---
$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.0303049096152694
$ws.Range("B2").Value = -0.0072368232351397
$ws.Range("A3").Value = -0.0072368232351397
$ws.Range("B3").Value = 0.0224459161341779

# --- sheet: gompertz cov ---
$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.0654227128172429
$ws.Range("B2").Value = -0.00711036413780034
$ws.Range("A3").Value = -0.00711036413780034
$ws.Range("B3").Value = 0.00164010407575795
